$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D/E retain plain-text storage so values such as
# "0.710" or "438.80" are not coerced into numbers (which would drop
# the significant trailing zero), matching the source data feed's
# text-formatted price/volume strings.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.104.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.154.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.01%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.146.68"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.07"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.666.57"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.127.50"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.146.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "489.26"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.72"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.710"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.66"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.42"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.36"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.24"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.02"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.27"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.05"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.74"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0747"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.54%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "438.80"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.120"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.914.86"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.260"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.03"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.44"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.26%  "
